# Calculated Fields and Data Hierarchies
# Update Sales (column D) values on Sheet1 and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Sales column values (D2:D10)
$ws.Range("D2").Value = 150000
$ws.Range("D3").Value = 200000
$ws.Range("D4").Value = 250000
$ws.Range("D5").Value = 850000
$ws.Range("D6").Value = 900000
$ws.Range("D7").Value = 110000
$ws.Range("D8").Value = 120000
$ws.Range("D9").Value = 60000
$ws.Range("D10").Value = 180000

# Move the active selection to H5 to match the saved view state
$ws.Activate()
$ws.Range("H5").Select()
